# Applies "minor fixing unit of measure" changes plus the small
# re-computed values that came along with the fix.
# Note: the PS parser used here does not accept scientific-notation
# numeric literals (e.g. 1.23E-2), so all values are written out in
# plain decimal form.

$wb = $excel.ActiveWorkbook

$wsSize = $wb.Worksheets.Item("Size")
$wsCost = $wb.Worksheets.Item("Cost")
$wsInd  = $wb.Worksheets.Item("Indicators")

# ---------------------------------------------------------------
# Sheet "Size"
# ---------------------------------------------------------------
$wsSize.Range("G2").Value = 181.613

$wsSize.Range("C3").Value = 0.24498372462774901
$wsSize.Range("D3").Value = 1.069245166222947
$wsSize.Range("E3").Value = 0.054720253761115452
$wsSize.Range("F3").Value = 0.047441902205855001
$wsSize.Range("G3").Value = 1.4163910468176659

# Columns C:G now share the same (wider) width, matching the new
# best-fit after G2's value grew from 0.181613 to 181.613.
# (ColumnWidth of 12 renders to a stored width of 12 + 5/6; subtract
# that offset so the stored OOXML width attribute comes out to 12.)
$wsSize.Range("C1:G1").EntireColumn.ColumnWidth = 11.166666666666666

# ---------------------------------------------------------------
# Sheet "Cost"
# ---------------------------------------------------------------
$wsCost.Range("I2").Value = 11.727660594589549

$wsCost.Range("E4").Value = 0.024498372462774899
$wsCost.Range("F4").Value = 0.1069245166222946
$wsCost.Range("G4").Value = 0.0054720253761115457
$wsCost.Range("H4").Value = 0.0047441902205855001
$wsCost.Range("I4").Value = 0.1416391046817666

$wsCost.Range("E6").Value = 0.0051294465965659201
$wsCost.Range("F6").Value = 0.022387756521829019
$wsCost.Range("G6").Value = 0.001145727618619054
$wsCost.Range("H6").Value = 0.0009933341660724839
$wsCost.Range("I6").Value = 0.029656264903086479

$wsCost.Range("I7").Value = 3.7650259026486861

$wsCost.Range("E8").Value = 1.9334530857065171
$wsCost.Range("F8").Value = 5.667545708009242
$wsCost.Range("G8").Value = 0.082659778913065252
$wsCost.Range("H8").Value = 0.02065686682046405
$wsCost.Range("I8").Value = 7.7043154394492879

# ---------------------------------------------------------------
# Sheet "Indicators"
# ---------------------------------------------------------------
$wsInd.Range("C2").Value = 159247.38033333331
$wsInd.Range("D2").Value = 167221.48974450849
$wsInd.Range("E2").Value = 326468.8700778418

$wsInd.Range("C5").Value = 0.00029999999999999997
$wsInd.Range("D5").Value = 0.00087878537564473449
$wsInd.Range("E5").Value = 0.00059646119930077983

$wsInd.Range("E6").Value = 0.25887889438258183
